$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Addr="D2"; Val="64.746.45"},
    @{Addr="E2"; Val="  -0.21%  "},
    @{Addr="D3"; Val="3.162.77"},
    @{Addr="E3"; Val="  -0.54%  "},
    @{Addr="E4"; Val="  +0.03%  "},
    @{Addr="D5"; Val="614.02"},
    @{Addr="E5"; Val="  +1.81%  "},
    @{Addr="D6"; Val="148.45"},
    @{Addr="E6"; Val="  -2.60%  "},
    @{Addr="E7"; Val="  -0.10%  "},
    @{Addr="D8"; Val="3.153.56"},
    @{Addr="E8"; Val="  -0.71%  "},
    @{Addr="D9"; Val="0.529"},
    @{Addr="E9"; Val="  -0.93%  "},
    @{Addr="E10"; Val="  -1.71%  "},
    @{Addr="D11"; Val="5.50"},
    @{Addr="E11"; Val="  -2.03%  "},
    @{Addr="D12"; Val="0.472"},
    @{Addr="E12"; Val="  -1.37%  "},
    @{Addr="D13"; Val="0.0000259"},
    @{Addr="E13"; Val="  -1.11%  "},
    @{Addr="D14"; Val="35.88"},
    @{Addr="E14"; Val="  -3.70%  "},
    @{Addr="D15"; Val="3.669.50"},
    @{Addr="E15"; Val="  -0.36%  "},
    @{Addr="E16"; Val="  +2.67%  "},
    @{Addr="D17"; Val="64.641.39"},
    @{Addr="E17"; Val="  -0.41%  "},
    @{Addr="D18"; Val="3.156.96"},
    @{Addr="E18"; Val="  +0.00%  "},
    @{Addr="D19"; Val="6.95"},
    @{Addr="E19"; Val="  -2.24%  "},
    @{Addr="D20"; Val="484.47"},
    @{Addr="E20"; Val="  -0.23%  "},
    @{Addr="D21"; Val="14.75"},
    @{Addr="E21"; Val="  -1.29%  "},
    @{Addr="D22"; Val="0.720"},
    @{Addr="E22"; Val="  +0.06%  "},
    @{Addr="D23"; Val="8.05"},
    @{Addr="E23"; Val="  +2.54%  "},
    @{Addr="D24"; Val="13.78"},
    @{Addr="E24"; Val="  -1.86%  "},
    @{Addr="D25"; Val="84.30"},
    @{Addr="E25"; Val="  -1.25%  "},
    @{Addr="E26"; Val="  +0.06%  "},
    @{Addr="D27"; Val="2.84"},
    @{Addr="E27"; Val="  -3.64%  "},
    @{Addr="D28"; Val="8.59"},
    @{Addr="E28"; Val="  -2.41%  "},
    @{Addr="D29"; Val="7.06"},
    @{Addr="E29"; Val="  -2.53%  "},
    @{Addr="E30"; Val="  -3.04%  "},
    @{Addr="D31"; Val="2.10"},
    @{Addr="E31"; Val="  -7.69%  "},
    @{Addr="E32"; Val="  -0.29%  "},
    @{Addr="E33"; Val="  +0.02%  "},
    @{Addr="D34"; Val="26.56"},
    @{Addr="E34"; Val="  -1.87%  "},
    @{Addr="D35"; Val="1.14"},
    @{Addr="E35"; Val="  +2.05%  "},
    @{Addr="D36"; Val="0.0₃0790"},
    @{Addr="E36"; Val="  +5.25%  "},
    @{Addr="D37"; Val="6.03"},
    @{Addr="E37"; Val="  -2.27%  "},
    @{Addr="D38"; Val="3.26"},
    @{Addr="E38"; Val="  -0.58%  "},
    @{Addr="D39"; Val="53.28"},
    @{Addr="E39"; Val="  -2.57%  "},
    @{Addr="D40"; Val="460.41"},
    @{Addr="E40"; Val="  +0.29%  "},
    @{Addr="D41"; Val="0.0402"},
    @{Addr="E41"; Val="  -1.14%  "},
    @{Addr="D42"; Val="0.121"},
    @{Addr="E42"; Val="  -5.80%  "},
    @{Addr="D43"; Val="8.42"},
    @{Addr="E43"; Val="  -1.81%  "},
    @{Addr="D44"; Val="2.853.78"},
    @{Addr="E44"; Val="  -2.14%  "},
    @{Addr="D45"; Val="2.34"},
    @{Addr="E45"; Val="  -5.37%  "},
    @{Addr="D46"; Val="0.270"},
    @{Addr="E46"; Val="  -2.94%  "},
    @{Addr="D47"; Val="2.49"},
    @{Addr="E47"; Val="  +5.16%  "},
    @{Addr="D48"; Val="26.66"},
    @{Addr="E48"; Val="  -2.40%  "},
    @{Addr="D50"; Val="0.115"},
    @{Addr="E50"; Val="  -1.68%  "},
    @{Addr="D51"; Val="120.55"},
    @{Addr="E51"; Val="  -0.03%  "}
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Addr)
    $c.NumberFormat = "@"
    $c.Value = $u.Val
    $c.Style = "Normal"
}
